# Revert "adding term 2.0.0"
$wb = $excel.ActiveWorkbook

# --- Metadata sheet: restore previous Version / Date / Contact values ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B3").Value = "1.1.0"
$meta.Range("B8").Value = "2023-07-10T23:08:03+02:00"
$meta.Range("B10").Value = "No display for ContactDetail"

# --- "Include from FSIII" sheet: drop the two concept rows that were added ---
$inc = $wb.Worksheets.Item("Include from FSIII")
$inc.Range("A2:B3").EntireRow.Delete()
